# Adding github account links
# Update the "Project website" row on the DataEntry sheet to become a
# "Project repository" row that points at the project's GitHub repo URL,
# and move the active cell selection from B4 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

$ws.Range("B4").Value = "https://github.com/publicusername/SensorEffector"
$ws.Range("A4").Value = "Project repository"

$ws.Activate()
$ws.Range("A4").Select()
